# NSMB workbook update - "Most of last stage done!"
# Adds new timing rows (203-216) with B/D data, inserts a couple of new
# checkpoint rows (207, 209) with notes in column G, and extends the
# shared D-column formula down through row 216.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Row 203: B/D data added; A/C already correct ---
$ws.Range("B203").Value = 69906
$ws.Range("D203").Formula = "=IF(B203 >  0,C203-B203, 0)"

# --- Row 204: B/D data added; A/C already correct ---
$ws.Range("B204").Value = 70786
$ws.Range("D204").Formula = "=IF(B204 >  0,C204-B204, 0)"

# --- Row 205: B/D data added; A/C already correct ---
$ws.Range("B205").Value = 71017
$ws.Range("D205").Formula = "=IF(B205 >  0,C205-B205, 0)"

# --- Row 206: B/D data added; A/C already correct ---
$ws.Range("B206").Value = 71467
$ws.Range("D206").Formula = "=IF(B206 >  0,C206-B206, 0)"

# --- Row 207: brand-new checkpoint row inserted before the old row 207 ---
$ws.Range("A207").Value = "Touch button"
$ws.Range("B207").Value = 71725
$ws.Range("C207").Value = 83509
$ws.Range("D207").Formula = "=IF(B207 >  0,C207-B207, 0)"
$ws.Range("G207").Value = "Note: not an absolute measure point"

# --- Row 208: former row 207 content, now pushed down, with B/D added ---
$ws.Range("A208").Value = "Enter door"
$ws.Range("B208").Value = 72016
$ws.Range("C208").Value = 83800
$ws.Range("D208").Formula = "=IF(B208 >  0,C208-B208, 0)"

# --- Row 209: brand-new checkpoint row ---
$ws.Range("A209").Value = "Checkpoint 1657"
$ws.Range("B209").Value = 72323
$ws.Range("C209").Value = 84128
$ws.Range("D209").Formula = "=IF(B209 >  0,C209-B209, 0)"
$ws.Range("G209").Value = "Approx (camera angle diffs)"

# --- Row 210 ---
$ws.Range("A210").Value = "Enter door"
$ws.Range("B210").Value = 72665
$ws.Range("C210").Value = 84487
$ws.Range("D210").Formula = "=IF(B210 >  0,C210-B210, 0)"

# --- Row 211 ---
$ws.Range("A211").Value = "Enter door"
$ws.Range("B211").Value = 73412
$ws.Range("C211").Value = 85239
$ws.Range("D211").Formula = "=IF(B211 >  0,C211-B211, 0)"

# --- Row 212 ---
$ws.Range("A212").Value = "Enter door"
$ws.Range("B212").Value = 73958
$ws.Range("C212").Value = 85786
$ws.Range("D212").Formula = "=IF(B212 >  0,C212-B212, 0)"

# --- Row 213 (no A value) ---
$ws.Range("B213").Value = 74154
$ws.Range("C213").Value = 85983
$ws.Range("D213").Formula = "=IF(B213 >  0,C213-B213, 0)"

# --- Row 214 (no A value) ---
$ws.Range("B214").Value = 74292
$ws.Range("C214").Value = 86121
$ws.Range("D214").Formula = "=IF(B214 >  0,C214-B214, 0)"

# --- Row 215 (no A value) ---
$ws.Range("B215").Value = 74361
$ws.Range("C215").Value = 86190
$ws.Range("D215").Formula = "=IF(B215 >  0,C215-B215, 0)"

# --- Row 216 (no A value) ---
$ws.Range("B216").Value = 74597
$ws.Range("C216").Value = 86424
$ws.Range("D216").Formula = "=IF(B216 >  0,C216-B216, 0)"

# Apply the same formatting (style index 16 - thin right border) that all
# the sibling data cells use, so the new cells render consistently (matches
# column A/B/C/D across the table).
$ws.Range("A207:D212").Borders.Item(10).LineStyle = 1
$ws.Range("B213:D216").Borders.Item(10).LineStyle = 1

# --- Update the frozen-pane scroll position / selection to match the
#     author's final view (pane frozen at row 1, scrolled so row 200 is
#     the first visible row beneath the freeze, with B217 selected). ---
$excel.ActiveWindow.ScrollRow = 200
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B217").Select()
